$wb = $excel.ActiveWorkbook

# Sheet1 ("names" sheet): remove the row-2 entry, shifting all subsequent
# ids up by one row (dimension shrinks from A1:A453 to A1:A452).
$names = $wb.Worksheets.Item("Sheet1")
$names.Rows.Item(2).Delete()

# Sheet2 ("used" sheet): append a new usage record in row 47
# (dimension grows from A1:C46 to A1:C47).
$used = $wb.Worksheets.Item("used")
$used.Range("A47").Value = "rhe1j7fa"
$used.Range("B47").Value = "ChatGPT Image 2026年1月24日 00_55_52.png"
$used.Range("C47").Value = "2026-01-24 00:58:06"
